# Auto-generated edit script: update cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.878.36'
$ws.Range("E2").Value = '  +4.43%  '

$ws.Range("D3").Value = '2.233.22'
$ws.Range("E3").Value = '  +3.27%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = "'259.80"
$ws.Range("E5").Value = '  +2.84%  '

$ws.Range("D6").Value = "'81.91"
$ws.Range("E6").Value = '  +11.63%  '

$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = '  +3.12%  '

$ws.Range("E8").Value = '  -0.25%  '

$ws.Range("D9").Value = "'0.604"
$ws.Range("E9").Value = '  +3.49%  '

$ws.Range("D10").Value = "'43.64"
$ws.Range("E10").Value = '  +9.85%  '

$ws.Range("D11").Value = "'0.0928"
$ws.Range("E11").Value = '  +2.51%  '

$ws.Range("D12").Value = "'7.05"
$ws.Range("E12").Value = '  +4.89%  '

$ws.Range("E13").Value = '  +2.54%  '

$ws.Range("D14").Value = '2.564.15'
$ws.Range("E14").Value = '  +3.05%  '

$ws.Range("D15").Value = "'14.65"
$ws.Range("E15").Value = '  +3.13%  '

$ws.Range("D16").Value = '2.237.74'
$ws.Range("E16").Value = '  +4.01%  '

$ws.Range("D17").Value = "'0.788"
$ws.Range("E17").Value = '  +2.59%  '

$ws.Range("D18").Value = '43.772.50'
$ws.Range("E18").Value = '  +4.43%  '

$ws.Range("E19").Value = '  +3.17%  '

$ws.Range("D20").Value = "'71.18"
$ws.Range("E20").Value = '  +0.90%  '

$ws.Range("D21").Value = "'6.05"
$ws.Range("E21").Value = '  +3.99%  '

$ws.Range("D22").Value = "'2.37"
$ws.Range("E22").Value = '  +9.86%  '

$ws.Range("D23").Value = "'232.64"
$ws.Range("E23").Value = '  +3.13%  '

$ws.Range("D24").Value = "'9.27"
$ws.Range("E24").Value = '  -2.58%  '

$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("D26").Value = "'10.81"
$ws.Range("E26").Value = '  +3.29%  '

$ws.Range("D27").Value = "'41.24"
$ws.Range("E27").Value = '  +12.37%  '

$ws.Range("E28").Value = '  +1.35%  '

$ws.Range("D29").Value = "'2.26"
$ws.Range("E29").Value = '  +3.35%  '

$ws.Range("E30").Value = '  -0.22%  '

$ws.Range("D31").Value = "'172.64"
$ws.Range("E31").Value = '  +2.35%  '

$ws.Range("D32").Value = "'0.0901"
$ws.Range("E32").Value = '  +13.53%  '

$ws.Range("D33").Value = "'20.64"
$ws.Range("E33").Value = '  +4.02%  '

$ws.Range("D34").Value = "'5.33"
$ws.Range("E34").Value = '  +5.09%  '

$ws.Range("E35").Value = '  +7.61%  '

$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").Value = "'0.123"
$ws.Range("E36").Value = '  +3.00%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.0371"
$ws.Range("E37").Value = '  +14.14%  '

$ws.Range("D38").Value = "'4.60"
$ws.Range("E38").Value = '  +8.55%  '

$ws.Range("D39").Value = "'13.03"
$ws.Range("E39").Value = '  +9.04%  '

$ws.Range("D40").Value = "'3.00"
$ws.Range("E40").Value = '  +26.63%  '

$ws.Range("D41").Value = "'2.14"
$ws.Range("E41").Value = '  +4.13%  '

$ws.Range("D42").Value = "'63.32"
$ws.Range("E42").Value = '  +8.37%  '

$ws.Range("D43").Value = "'5.53"
$ws.Range("E43").Value = '  +7.91%  '

$ws.Range("D44").Value = "'0.203"
$ws.Range("E44").Value = '  +4.08%  '

$ws.Range("D45").Value = "'104.29"
$ws.Range("E45").Value = '  +4.24%  '

$ws.Range("E46").Value = '  +2.44%  '

$ws.Range("D47").Value = "'0.0985"
$ws.Range("E47").Value = '  +2.63%  '

$ws.Range("E48").Value = '  +4.03%  '

$ws.Range("E49").Value = '  +28.86%  '

$ws.Range("D50").Value = "'0.442"
$ws.Range("E50").Value = '  -5.28%  '

$ws.Range("E51").Value = '  +3.71%  '
